$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.128047585487366
$ws.Range("B1").Value = 2.279197216033936
$ws.Range("C1").Value = 10.95813369750977
$ws.Range("D1").Value = 2.010145902633667
$ws.Range("E1").Value = 1.283152103424072
